$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 131
$ws1.Range("F6").Value = 20
$ws1.Range("F8").Value = 4802
$ws1.Range("F9").Value = 4802
$ws1.Range("F11").Value = 129
$ws1.Range("F14").Value = 1100
$ws1.Range("F15").Value = 632
$ws1.Range("F16").Value = 4401
$ws1.Range("F17").Value = 173
$ws1.Range("F18").Value = 174
$ws1.Range("F19").Value = 76
$ws1.Range("F20").Value = 225
$ws1.Range("F21").Value = 3529
$ws1.Range("F25").Value = 3199
$ws1.Range("F27").Value = 132
$ws1.Range("F32").Value = 85
$ws1.Range("F33").Value = 66
$ws1.Range("F37").Value = 5596
$ws1.Range("F38").Value = 864
$ws1.Range("F39").Value = 412
$ws1.Range("F43").Value = 1135
$ws1.Range("F44").Value = 507
$ws1.Range("F46").Value = 2015
$ws1.Range("F47").Value = 300
$ws1.Range("F48").Value = 69
$ws1.Range("F49").Value = 708
$ws1.Range("F50").Value = 859

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 66

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 131
$ws4.Range("F7").Value = 20
$ws4.Range("F9").Value = 4802
$ws4.Range("F10").Value = 4802
$ws4.Range("F13").Value = 66
$ws4.Range("F16").Value = 1100
$ws4.Range("F17").Value = 632
$ws4.Range("F18").Value = 4401
$ws4.Range("F19").Value = 173
$ws4.Range("F20").Value = 174
$ws4.Range("F21").Value = 76
$ws4.Range("F22").Value = 225
$ws4.Range("F23").Value = 3529
$ws4.Range("F24").Value = 3200
$ws4.Range("F26").Value = 132
$ws4.Range("F30").Value = 85
$ws4.Range("F31").Value = 66
$ws4.Range("F36").Value = 5596
$ws4.Range("F38").Value = 864
$ws4.Range("F39").Value = 412
$ws4.Range("F45").Value = 1135
$ws4.Range("F46").Value = 507
$ws4.Range("F47").Value = 2015
$ws4.Range("F48").Value = 300
$ws4.Range("F49").Value = 859
